$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Freesia"
$ws.Range("B3").Value = "Freesia1"

$ws.Range("A4").Value = "a_guy"
$ws.Range("B4").Value = "hhh"

$ws.Range("A5").Value = "dd"
$ws.Range("B5").Value = "ff"

$ws.Range("A6").Value = "iii"
$ws.Range("B6").Value = "iii"

$ws.Range("A7").Value = "seemore"
$ws.Range("B7").Value = "butts"
